{"js": "// Update the worksheet date line and regenerate all 100 addition/subtraction\n// answers in the 20x5 table with a new set of equations (new \"answer key\"\n// output for commit c8c62b6).\n\nconst NEW_DATE = \"2026-02-13 Friday\";\n\n// Row-major (top-to-bottom, left-to-right) replacement values for every\n// cell of the 20x5 table, in the exact order the cells appear in the\n// document.\nconst NEW_ANSWERS = [\n  \"86-67=19\", \"52+16=68\", \"28+51=79\", \"70+22=92\", \"19-5=14\",\n  \"88-20=68\", \"13-6=7\", \"62-5=57\", \"36+7=43\", \"86-16=70\",\n  \"51-3=48\", \"74-52=22\", \"41+55=96\", \"54-45=9\", \"12-4=8\",\n  \"45+35=80\", \"8+15=23\", \"72-62=10\", \"78-4=74\", \"93-65=28\",\n  \"14-11=3\", \"27+66=93\", \"28+26=54\", \"61+6=67\", \"20+57=77\",\n  \"18+43=61\", \"40-29=11\", \"23-17=6\", \"64-10=54\", \"78+14=92\",\n  \"67+28=95\", \"32-5=27\", \"2+13=15\", \"86-62=24\", \"86-13=73\",\n  \"76-36=40\", \"93-22=71\", \"17+6=23\", \"62-4=58\", \"9+60=69\",\n  \"3+53=56\", \"62-0=62\", \"58-29=29\", \"85-69=16\", \"38+6=44\",\n  \"85-72=13\", \"91-85=6\", \"65+27=92\", \"30-2=28\", \"96+1=97\",\n  \"95-70=25\", \"33+17=50\", \"31+38=69\", \"13+51=64\", \"92+0=92\",\n  \"93-42=51\", \"87-39=48\", \"50+15=65\", \"69-9=60\", \"24+2=26\",\n  \"46-28=18\", \"96-9=87\", \"95-32=63\", \"70+8=78\", \"68-18=50\",\n  \"91-87=4\", \"75-25=50\", \"44+34=78\", \"61+3=64\", \"42-39=3\",\n  \"67-43=24\", \"24+11=35\", \"80-29=51\", \"23-21=2\", \"51+46=97\",\n  \"39+58=97\", \"52-42=10\", \"38+56=94\", \"85+3=88\", \"43+11=54\",\n  \"30+29=59\", \"77-8=69\", \"65+5=70\", \"85-22=63\", \"65-45=20\",\n  \"65+33=98\", \"22+39=61\", \"1+66=67\", \"1+21=22\", \"20+3=23\",\n  \"66+23=89\", \"21+11=32\", \"97-38=59\", \"74-42=32\", \"89-29=60\",\n  \"1+42=43\", \"59-15=44\", \"26+6=32\", \"51-19=32\", \"90-75=15\",\n];\n\n// 1) Update the date paragraph at the top of the document.\nconst firstPara = context.document.body.paragraphs.getFirst();\nfirstPara.insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// 2) Update every answer cell in the (only) table, row by row, left to\n// right, matching the order of NEW_ANSWERS above.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\nfor (const row of rows.items) {\n  row.load(\"cellCount\");\n}\nawait context.sync();\n\nlet k = 0;\nfor (let r = 0; r < rows.items.length; r++) {\n  const colCount = rows.items[r].cellCount;\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = NEW_ANSWERS[k];\n    k++;\n  }\n}\n\nif (k !== NEW_ANSWERS.length) {\n  throw new Error(\n    `Expected ${NEW_ANSWERS.length} table cells but updated ${k}.`\n  );\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date line and regenerate all 100 addition/subtraction\n# answers in the 20x5 table with a new set of equations (new \"answer key\"\n# output for commit c8c62b6).\n\n$d = $word.ActiveDocument\n\n$NewDate = \"2026-02-13 Friday\"\n\n# Row-major (top-to-bottom, left-to-right) replacement values for every\n# cell of the 20x5 table, in the exact order the cells appear in the\n# document.\n$NewAnswers = @(\n    \"86-67=19\", \"52+16=68\", \"28+51=79\", \"70+22=92\", \"19-5=14\",\n    \"88-20=68\", \"13-6=7\", \"62-5=57\", \"36+7=43\", \"86-16=70\",\n    \"51-3=48\", \"74-52=22\", \"41+55=96\", \"54-45=9\", \"12-4=8\",\n    \"45+35=80\", \"8+15=23\", \"72-62=10\", \"78-4=74\", \"93-65=28\",\n    \"14-11=3\", \"27+66=93\", \"28+26=54\", \"61+6=67\", \"20+57=77\",\n    \"18+43=61\", \"40-29=11\", \"23-17=6\", \"64-10=54\", \"78+14=92\",\n    \"67+28=95\", \"32-5=27\", \"2+13=15\", \"86-62=24\", \"86-13=73\",\n    \"76-36=40\", \"93-22=71\", \"17+6=23\", \"62-4=58\", \"9+60=69\",\n    \"3+53=56\", \"62-0=62\", \"58-29=29\", \"85-69=16\", \"38+6=44\",\n    \"85-72=13\", \"91-85=6\", \"65+27=92\", \"30-2=28\", \"96+1=97\",\n    \"95-70=25\", \"33+17=50\", \"31+38=69\", \"13+51=64\", \"92+0=92\",\n    \"93-42=51\", \"87-39=48\", \"50+15=65\", \"69-9=60\", \"24+2=26\",\n    \"46-28=18\", \"96-9=87\", \"95-32=63\", \"70+8=78\", \"68-18=50\",\n    \"91-87=4\", \"75-25=50\", \"44+34=78\", \"61+3=64\", \"42-39=3\",\n    \"67-43=24\", \"24+11=35\", \"80-29=51\", \"23-21=2\", \"51+46=97\",\n    \"39+58=97\", \"52-42=10\", \"38+56=94\", \"85+3=88\", \"43+11=54\",\n    \"30+29=59\", \"77-8=69\", \"65+5=70\", \"85-22=63\", \"65-45=20\",\n    \"65+33=98\", \"22+39=61\", \"1+66=67\", \"1+21=22\", \"20+3=23\",\n    \"66+23=89\", \"21+11=32\", \"97-38=59\", \"74-42=32\", \"89-29=60\",\n    \"1+42=43\", \"59-15=44\", \"26+6=32\", \"51-19=32\", \"90-75=15\"\n)\n\n# 1) Update the date paragraph at the top of the document.\n$d.Paragraphs.First.Range.Text = $NewDate\n\n# 2) Update every answer cell in the (only) table, row by row, left to\n# right, matching the order of $NewAnswers above.\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$k = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $NewAnswers[$k]\n        $k++\n    }\n}\n\nif ($k -ne $NewAnswers.Count) {\n    throw \"Expected $($NewAnswers.Count) table cells but updated $k.\"\n}\n"}
